$wb = $excel.ActiveWorkbook

# --- Clear selection state on Italy first (it currently is B4 / tabSelected) ---
$italy = $wb.Worksheets.Item("Italy")

# Duplicate the "Italy" sheet (keeps formatting/merges/styles) and place it right after Italy.
$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# --- Update the market-specific cell contents on the new Spain sheet ---
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2035 "

# --- Column widths to match the Spain sheet's layout ---
$spain.Columns("A:A").ColumnWidth = 25.109375
$spain.Columns("B:B").ColumnWidth = 21.6640625
$spain.Columns("C:C").ColumnWidth = 16.44140625
$spain.Columns("D:D").ColumnWidth = 25.33203125

# --- Row heights for the (now taller) header rows ---
$spain.Rows("3:3").RowHeight = 28.8
$spain.Rows("4:4").RowHeight = 28.8
$spain.Rows("5:5").RowHeight = 28.8

# --- Page setup (portrait) on both the copied Italy sheet and new Spain sheet ---
$italy.PageSetup.Orientation = 1
$spain.PageSetup.Orientation = 1

# --- Selection / active cell bookkeeping ---
$italy.Range("A1:D13").Select()
$spain.Range("D9").Select()

# Spain becomes the active/visible tab, as in the authored workbook.
$spain.Activate()
